$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: F1/H1 text shifts (C1->C5, C2->C6) as the "atom map" columns
# are renumbered to make room for the two new side-chain columns.
$ws.Range("F1").Value = "C5"
$ws.Range("H1").Value = "C6"

# New header cells for the two appended side-alkyl-chain atom columns.
$ws.Range("I1").Value = "C1"
$ws.Range("J1").Value = "C2"

# Match the bold/centered/bordered header formatting already used by B1:H1.
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data rows: new I/J columns with the matched side-chain atom names.
$ws.Range("I2").Value = "C5"
$ws.Range("J2").Value = "C4"
$ws.Range("I3").Value = "C5"
$ws.Range("J3").Value = "C4"
$ws.Range("I4").Value = "C5"
$ws.Range("J4").Value = "C4"
$ws.Range("I5").Value = "C5"
$ws.Range("J5").Value = "C4"
$ws.Range("I6").Value = "C8"
$ws.Range("J6").Value = "C7"
$ws.Range("I7").Value = "C8"
$ws.Range("J7").Value = "C7"
$ws.Range("I8").Value = "C6"
$ws.Range("J8").Value = "C5"
$ws.Range("I9").Value = "C6"
$ws.Range("J9").Value = "C5"
$ws.Range("I10").Value = "C9"
$ws.Range("J10").Value = "C8"
$ws.Range("I11").Value = "C7"
$ws.Range("J11").Value = "C6"
$ws.Range("I12").Value = "C7"
$ws.Range("J12").Value = "C6"
$ws.Range("I13").Value = "C1"
$ws.Range("J13").Value = "C2"
$ws.Range("I14").Value = "C6"
$ws.Range("J14").Value = "C5"
